# Sales_and_marketing/data.xlsx -- "Add files via upload"
#
# The attendance sheet used four different spellings for the same two
# statuses ("Present"/"present" and "ABSENT"/"Absent"). This normalises
# every attendance cell (columns D:P, rows 2:44) down to the single-letter
# codes "P" and "A", and updates the active sheet/selection bookmarks left
# behind by the author's last save (attendance tab became the active tab,
# with G11 selected; the Introduction tab kept its own selection but lost
# "active" status).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attendance")

for ($r = 2; $r -le 44; $r++) {
    for ($c = 4; $c -le 16; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Text
        if ($v -eq "Present" -or $v -eq "present") {
            $cell.Value = "P"
        } elseif ($v -eq "ABSENT" -or $v -eq "Absent") {
            $cell.Value = "A"
        }
    }
}

# Make "attendance" the active sheet, with the new selection left on it,
# and restore the "Introduction to sales and marke" sheet's own selection
# (no longer the active tab).
$ws3 = $wb.Worksheets.Item("Introduction to sales and marke")
$ws3.Range("I21").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("G11").Select() | Out-Null
